# The match-innings rows for Vijay Shankar (Sunrisers Hyderabad) got
# reordered/updated as more data became available ("updated activity till
# excel form"). Row 3 (52/51/6/0) is untouched. Rows 2, 4, 5, 6 get new
# runs/balls/fours/sixes values for columns C:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (runs, balls, fours, sixes)
$newValues = @{
    2 = @("26", "27", "4", "0")
    4 = @("12", "7", "0", "1")
    5 = @("7", "10", "0", "0")
    6 = @("0", "1", "0", "0")
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $rng = $ws.Range("C$row`:F$row")
    # Keep these as text (matching the existing "numbers stored as text"
    # formatting used throughout the sheet) rather than converting to
    # numeric values.
    $rng.NumberFormat = "@"
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("D$row").Value = $vals[1]
    $ws.Range("E$row").Value = $vals[2]
    $ws.Range("F$row").Value = $vals[3]
}
